# Commit: "Update Work Week and Social Spending"
#
# Refresh the Switzerland "GDP per Capita" series on the Data sheet:
#   - every existing year (1851-2010, rows 2-161) gets a revised value
#   - six newly available years (2011-2016, rows 162-167) are appended
#
# Column E holds these figures as TEXT in this workbook (not numbers), so the
# Text number format is applied before the values are written; otherwise
# Excel would auto-coerce the numeric-looking strings into real numbers.
# The Metadata sheet's citation text is untouched -- it does not change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# GDP-per-capita figures, one per year, in order from 1851 through 2016.
$gdpValues = @(
    "2444.22610239988","2445.85217664588","2379.2794706162","2263.80653910788","2372.84918976211","2419.27101505315","2785.79517964569","2861.66667495646","2894.26277033012","2933.63057110256","2948.31278589882","2822.51362876576","2879.02433588928","2807.79281136665","2792.1978335779","2814.7017182278","2746.02266573331","3053.11450938051","3041.02859434905","2954.37649900731",
    "3063.71744352942","3140.48686963731","3320.86142093715","3322.5070206781","3598.77225366935","3851.91415709744","3630.18782609982","3485.83696086339","3509.22223480297","3869.36628148129","3943.41551086938","3998.66157597033","4194.24861044122","4542.09409669602","4772.86416899563","5036.99075848235","5172.16849930844","5118.90121485149","5137.32807514365","5309.32048505504",
    "4937.06615271165","5270.18587639606","5339.8903028207","5337.31207346907","5632.35874744704","5967.24507057932","6241.89066472412","6427.50876279354","6337.31853386101","6611.77821356305","6398.89186862999","6616.154464283","6472.23998127461","6879.57492621256","7184.5796975677","7584.05690252643","7464.97944589993","7507.70959998377","7858.35886708257","8048.35153312609",
    "8043.94146001353","8311.78027545233","8172.06013176882","7595.28553094459","7459.25889916841","8079.41922593748","7840.90390719607","7115.60757090518","6747.26797745246","7363.9596412385","5999.14144316014","7181.95352888635","7670.23768949731","8195.52064929545","8385.2867756263","8656.33106890515","9259.93620280581","9658.81387249873","9920.03754250447","9969.3240645101",
    "10055.4851996385","9752.44296534974","10009.108560589","9997.37644988818","9478.80775609377","9301.84625825247","10039.2922320319","9608.57001962028","9836.15191455033","9802.37306877984","9217.34890451972","8739.76696861602","8642.97285805962","8890.73881147861","8839.31526323939","10866.098130433","11020.9505631674","11189.3430725298","10943.5440327929","11541.0375505428",
    "11877.89329403","12016.308232061","12470.9856694293","13171.4062803502","13999.4878426012","14507.0601800455","14667.650420408","14644.4582568088","15469.5914767422","16358.3594046842","17307.5627469558","17793.9764899977","18437.8389135498","19195.446462746","19802.512491383","20224.5381975754","20742.2525407236","21395.3863755189","22307.672926555","23478.991767023",
    "24486.9781291002","25356.9109020439","26003.7151749567","25718.2703729681","24467.991703049","24752.9073088838","25430.5674454033","26158.8002988737","26641.0949518161","27405.7302362483","27729.584967226","27660.9484973249","28009.3290841041","28315.0825944092","29416.75224461","30847.9942567071","31774.3711593065","32482.9076702734","33315.5867704057","34250",
    "34128.6419071025","34356.8020946259","34625.5336467752","35423.8278284248","36004.2298937464","36720.4904727663","38165.477352616","39890.5375751659","41091.6368773432","43251.1037536027","44398.6387814687","44934.3167460536","45441.5992803626","47263.5967262545","49273.9023593869","51862.9647705908","54509.5396598574","56054.8838155727","55171.5007076847","57219.4983120623",
    "58686","58674","59036","59485","59307","59662"
)

$firstYear = 1851
$firstRow = 2
$existingLastRow = 161
$countryCode = 756
$countryName = "Switzerland"
$indicatorName = "GDP per Capita"

$lastRow = $firstRow + $gdpValues.Count - 1

$ws.Range("E" + $firstRow + ":E" + $lastRow).NumberFormat = "@"

for ($i = 0; $i -lt $gdpValues.Count; $i++) {
    $row = $firstRow + $i
    $year = $firstYear + $i

    if ($row -gt $existingLastRow) {
        # Brand-new row (years 2011-2016) -- populate every column.
        $ws.Cells.Item($row, 1).Value = $countryCode
        $ws.Cells.Item($row, 2).Value = $countryName
        $ws.Cells.Item($row, 3).Value = $indicatorName
        $ws.Cells.Item($row, 4).Value = $year
    }

    $ws.Cells.Item($row, 5).Value = $gdpValues[$i]
}
